$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the (now-redundant) Network/H column values on a handful of
#     2023 rows (H32, H34, H35, H36, H37). Clear() drops the <c> element
#     entirely rather than leaving an empty styled cell behind. ---
$ws.Range("H32").Clear()
$ws.Range("H34").Clear()
$ws.Range("H35").Clear()
$ws.Range("H36").Clear()
$ws.Range("H37").Clear()

# --- Insert a new row for the "2023_TM160_IPA_29" model run right above
#     the old row 41 (2025 Plan run), pushing everything below down by one. ---
$ws.Rows.Item(41).EntireRow.Insert()

$ws.Range("A41").Value = 2023
$ws.Range("B41").Value = "2023_TM160_IPA_29"
$ws.Range("C41").Value = "RTP2025_IP"
$ws.Range("D41").Value = "Past year"
$ws.Range("E41").Value = "Refined empres estimates, WFH remains at ~25%"
$ws.Range("F41").Value = "petrale"
$ws.Range("G41").Value = "n/a"
$ws.Range("H41").Value = "current"
$ws.Range("I41").Value = "BlueprintNetworks_v10\net_2023_Blueprint"
$ws.Range("J41").Value = "model2-b"
$ws.Range("K41").Value = "https://app.asana.com/0/1204085012544660/1205838169476835/f"
$ws.Range("L41").Value = 17.77
$ws.Range("M41").Value = "na"
$ws.Range("N41").Value = "na"
$ws.Range("O41").Value = 0.94
$ws.Range("P41").Value = 0.855
$ws.Range("Q41").Value = 120
$ws.Range("R41").Value = 0
$ws.Range("S41").Value = 45

# --- Move the active selection to A41, matching the author's final cursor
#     position after adding the new row. ---
$ws.Range("A41").Select()
